$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A184").Value = "20-09-2021"
$ws.Range("B184").Value = 4.38
$ws.Range("C184").Value = 3.91
$ws.Range("D184").Value = 3.46
$ws.Range("E184").Value = 3.38

$ws.Range("A185").Value = "21-09-2021"
$ws.Range("B185").Value = 4.34
$ws.Range("C185").Value = 3.86
$ws.Range("D185").Value = 3.41
$ws.Range("E185").Value = 3.35

$ws.Range("A186").Value = "22-09-2021"
$ws.Range("B186").Value = 4.34
$ws.Range("C186").Value = 3.89
$ws.Range("D186").Value = 3.43
$ws.Range("E186").Value = 3.35

$ws.Range("A187").Value = "23-09-2021"
$ws.Range("B187").Value = 4.48
$ws.Range("C187").Value = 3.87
$ws.Range("D187").Value = 3.45
$ws.Range("E187").Value = 3.36

$ws.Range("A188").Value = "24-09-2021"
$ws.Range("B188").Value = 4.66
$ws.Range("C188").Value = 3.91
$ws.Range("D188").Value = 3.46
$ws.Range("E188").Value = 3.35
